$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "59.725.65"
$ws.Range("E2").Value = "  +0.73%  "
Set-TextValue "D3" "2.652.33"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").Value = "  -0.11%  "
Set-TextValue "D5" "537.87"
$ws.Range("E5").Value = "  -1.03%  "
Set-TextValue "D6" "146.92"
$ws.Range("E6").Value = "  +4.23%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +1.36%  "
Set-TextValue "D9" "6.85"
$ws.Range("E9").Value = "  +6.05%  "
Set-TextValue "D10" "0.103"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("E12").Value = "  +0.12%  "
Set-TextValue "D13" "3.128.32"
$ws.Range("E13").Value = "  +2.06%  "
Set-TextValue "D14" "59.631.10"
$ws.Range("E14").Value = "  +0.70%  "
Set-TextValue "D15" "21.47"
$ws.Range("E15").Value = "  +4.58%  "
Set-TextValue "D16" "2.662.37"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("E17").Value = "  +1.27%  "
Set-TextValue "D18" "4.49"
$ws.Range("E18").Value = "  +2.93%  "
Set-TextValue "D19" "340.81"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("E20").Value = "  +2.43%  "
Set-TextValue "D21" "6.23"
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("E22").Value = "  +0.07%  "
Set-TextValue "D23" "66.70"
$ws.Range("E23").Value = "  -1.19%  "
Set-TextValue "D24" "0.419"
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("E25").Value = "  -0.11%  "
Set-TextValue "D26" "0.997"
$ws.Range("E26").Value = "  -0.21%  "
Set-TextValue "D27" "7.33"
$ws.Range("E27").Value = "  +1.60%  "
Set-TextValue "D28" "0.0₃0752"
$ws.Range("E28").Value = "  +2.19%  "
Set-TextValue "D29" "0.998"
$ws.Range("E30").Value = "  -2.63%  "
Set-TextValue "D31" "5.89"
$ws.Range("E31").Value = "  +1.60%  "
Set-TextValue "D32" "18.91"
$ws.Range("E32").Value = "  +0.92%  "
Set-TextValue "D33" "150.74"
$ws.Range("E33").Value = "  +0.56%  "
Set-TextValue "D34" "4.02"
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("E36").Value = "  +3.57%  "
Set-TextValue "D37" "0.844"
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("E39").Value = "  +1.93%  "
Set-TextValue "D40" "286.71"
$ws.Range("E40").Value = "  +3.68%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  +1.89%  "
Set-TextValue "D43" "10.75"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  +3.05%  "
Set-TextValue "D45" "19.35"
$ws.Range("E45").Value = "  +3.89%  "
Set-TextValue "D46" "0.0949"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D47" "4.68"
$ws.Range("E47").Value = "  +3.88%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D48" "0.0228"
$ws.Range("E48").Value = "  +2.07%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D49" "1.967.53"
$ws.Range("E49").Value = "  +1.31%  "
Set-TextValue "D50" "18.47"
$ws.Range("E50").Value = "  +0.56%  "
Set-TextValue "D51" "112.34"
$ws.Range("E51").Value = "  +1.18%  "
